$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws1.Range("H18").Value = 2874.75
$ws1.Range("I18").Value = 2833
$ws1.Range("K18").Value = 2833
$ws1.Range("M18").Value = -2549
$ws1.Range("H64").Value = 5881
$ws1.Range("I64").Value = 4262.25
$ws1.Range("J64").Value = 7499.75
$ws1.Range("K64").Value = 4262.25
$ws1.Range("L64").Value = 7499.75
$ws1.Range("M64").Value = -4014.25
$ws1.Range("N64").Value = -7995.75
$ws1.Range("H67").Value = 5881
$ws1.Range("I67").Value = 4262.25
$ws1.Range("J67").Value = 7499.75
$ws1.Range("K67").Value = 4262.25
$ws1.Range("L67").Value = 7499.75
$ws1.Range("M67").Value = -3404.25
$ws1.Range("N67").Value = -9215.75
$ws1.Range("H87").Value = 97587.75
$ws1.Range("J87").Value = 97587.75
$ws1.Range("L87").Value = 97587.75
$ws1.Range("N87").Value = -100083.75
$ws1.Range("H90").Value = 97587.75
$ws1.Range("J90").Value = 97587.75
$ws1.Range("L90").Value = 292763.25
$ws1.Range("N90").Value = -305243.25
$ws1.Range("H137").Value = 5575.7896
$ws1.Range("I137").Value = 3599
$ws1.Range("J137").Value = 6281.7856
$ws1.Range("K137").Value = 10797
$ws1.Range("L137").Value = 18845.3568
$ws1.Range("M137").Value = -8247
$ws1.Range("N137").Value = -23945.3568

$ws2 = $wb.Worksheets.Item("ARM")
$ws2.Range("H11").Value = 3668.6667
$ws2.Range("I11").Value = 3668.6667
$ws2.Range("K11").Value = 3668.6667
$ws2.Range("M11").Value = -3524.6667
$ws2.Range("H12").Value = 1000
$ws2.Range("I12").Value = 0
$ws2.Range("K12").Value = 0
$ws2.Range("M12").ClearContents()
$ws2.Range("H14").Value = 258.42856
$ws2.Range("I14").Value = 102.25
$ws2.Range("J14").Value = 466.66666
$ws2.Range("K14").Value = 102.25
$ws2.Range("L14").Value = 466.66666
$ws2.Range("M14").Value = 72.75
$ws2.Range("N14").Value = -816.66666
$ws2.Range("H30").Value = 1175
$ws2.Range("I30").Value = 1175
$ws2.Range("J30").Value = 0
$ws2.Range("K30").Value = 1175
$ws2.Range("L30").Value = 0
$ws2.Range("M30").Value = -1025
$ws2.Range("N30").ClearContents()
$ws2.Range("H102").Value = 2013.8
$ws2.Range("I102").Value = 2013.8
$ws2.Range("J102").Value = 0
$ws2.Range("K102").Value = 2013.8
$ws2.Range("L102").Value = 0
$ws2.Range("M102").Value = -391.8
$ws2.Range("N102").ClearContents()

$ws3 = $wb.Worksheets.Item("BSM")
$ws3.Range("H37").Value = 500
$ws3.Range("I37").Value = 500
$ws3.Range("J37").Value = 0
$ws3.Range("K37").Value = 500
$ws3.Range("L37").Value = 0
$ws3.Range("M37").Value = -363
$ws3.Range("N37").ClearContents()
$ws3.Range("H105").Value = 2954.4736
$ws3.Range("I105").Value = 1700.1428
$ws3.Range("J105").Value = 3686.1667
$ws3.Range("K105").Value = 1700.1428
$ws3.Range("L105").Value = 3686.1667
$ws3.Range("M105").Value = 46.85719999999992
$ws3.Range("N105").Value = -7180.1667

$ws4 = $wb.Worksheets.Item("CRP")
$ws4.Range("H6").Value = 470
$ws4.Range("J6").Value = 0
$ws4.Range("L6").Value = 0
$ws4.Range("N6").ClearContents()
$ws4.Range("H11").Value = 9400
$ws4.Range("I11").Value = 10101
$ws4.Range("J11").Value = 9166.333000000001
$ws4.Range("K11").Value = 10101
$ws4.Range("L11").Value = 9166.333000000001
$ws4.Range("M11").Value = -9961
$ws4.Range("N11").Value = -9446.333000000001
$ws4.Range("H19").Value = 22333582
$ws4.Range("I19").Value = 22333582
$ws4.Range("K19").Value = 22333582
$ws4.Range("M19").Value = -22333412
$ws4.Range("H24").Value = 22333582
$ws4.Range("I24").Value = 22333582
$ws4.Range("K24").Value = 22333582
$ws4.Range("M24").Value = -22333412
$ws4.Range("H31").Value = 3925.0625
$ws4.Range("I31").Value = 3152.8462
$ws4.Range("K31").Value = 3152.8462
$ws4.Range("M31").Value = -2857.8462
$ws4.Range("H32").Value = 2627.5
$ws4.Range("I32").Value = 2627.5
$ws4.Range("J32").Value = 0
$ws4.Range("K32").Value = 2627.5
$ws4.Range("L32").Value = 0
$ws4.Range("M32").Value = -2311.5
$ws4.Range("N32").ClearContents()
$ws4.Range("H34").Value = 3925.0625
$ws4.Range("I34").Value = 3152.8462
$ws4.Range("K34").Value = 3152.8462
$ws4.Range("M34").Value = -2950.8462
$ws4.Range("H62").Value = 0
$ws4.Range("I62").Value = 0
$ws4.Range("K62").Value = 0
$ws4.Range("M62").ClearContents()
$ws4.Range("H65").Value = 0
$ws4.Range("I65").Value = 0
$ws4.Range("K65").Value = 0
$ws4.Range("M65").ClearContents()
$ws4.Range("H99").Value = 2479.6
$ws4.Range("I99").Value = 1500
$ws4.Range("K99").Value = 1500
$ws4.Range("M99").Value = -2
$ws4.Range("H126").Value = 2479.6
$ws4.Range("I126").Value = 1500
$ws4.Range("K126").Value = 4500
$ws4.Range("M126").Value = -2030

$ws5 = $wb.Worksheets.Item("CUL")
$ws5.Range("H6").Value = 60
$ws5.Range("I6").Value = 60
$ws5.Range("K6").Value = 180
$ws5.Range("M6").Value = -67
$ws5.Range("H26").Value = 1265.2325
$ws5.Range("I26").Value = 1155.5264
$ws5.Range("J26").Value = 2099
$ws5.Range("K26").Value = 3466.5792
$ws5.Range("L26").Value = 6297
$ws5.Range("M26").Value = -3178.5792
$ws5.Range("N26").Value = -6873
$ws5.Range("H47").Value = 54.75
$ws5.Range("I47").Value = 54.75
$ws5.Range("K47").Value = 164.25
$ws5.Range("M47").Value = 266.75
$ws5.Range("H107").Value = 460.58334
$ws5.Range("J107").Value = 485.25
$ws5.Range("L107").Value = 1455.75
$ws5.Range("N107").Value = -5295.75

$ws6 = $wb.Worksheets.Item("GSM")
$ws6.Range("H13").Value = 1093.6
$ws6.Range("I13").Value = 498.5
$ws6.Range("J13").Value = 1490.3334
$ws6.Range("K13").Value = 498.5
$ws6.Range("L13").Value = 1490.3334
$ws6.Range("M13").Value = -359.5
$ws6.Range("N13").Value = -1768.3334
$ws6.Range("H102").Value = 1907.8334
$ws6.Range("I102").Value = 1889.4
$ws6.Range("K102").Value = 1889.4
$ws6.Range("M102").Value = -267.4000000000001
$ws6.Range("H113").Value = 1266.1538
$ws6.Range("J113").Value = 1674
$ws6.Range("L113").Value = 1674
$ws6.Range("N113").Value = -6014
$ws6.Range("H126").Value = 0
$ws6.Range("I126").Value = 0
$ws6.Range("K126").Value = 0
$ws6.Range("M126").ClearContents()

$ws7 = $wb.Worksheets.Item("LTW")
$ws7.Range("H9").Value = 306
$ws7.Range("I9").Value = 306
$ws7.Range("K9").Value = 306
$ws7.Range("M9").Value = -82
$ws7.Range("H31").Value = 1935
$ws7.Range("I31").Value = 1322
$ws7.Range("K31").Value = 1322
$ws7.Range("M31").Value = -1074

$ws8 = $wb.Worksheets.Item("WVR")
$ws8.Range("H19").Value = 11450
$ws8.Range("I19").Value = 15900
$ws8.Range("J19").Value = 7000
$ws8.Range("K19").Value = 15900
$ws8.Range("L19").Value = 7000
$ws8.Range("M19").Value = -15726
$ws8.Range("N19").Value = -7348
$ws8.Range("H136").Value = 5857.7144
$ws8.Range("I136").Value = 5334
$ws8.Range("J136").Value = 9000
$ws8.Range("K136").Value = 16002
$ws8.Range("L136").Value = 27000
$ws8.Range("M136").Value = -13452
$ws8.Range("N136").Value = -32100
